# Update "IEnumerable vs IQueryable" worksheet:
#   - Row 6 (Where example) now calls .AsQueryable()/.AsEnumerable() BEFORE
#     .Where(...) instead of after it, and the matching SQL loses the stray
#     leading whitespace and the (now redundant) WHERE clause in column D.
#   - Row 7/8 (Select example) similarly moves .AsEnumerable()/.AsQueryable()
#     before .Select(...), and columns D7/D8 are corrected to the
#     s.Id/MaintainerEmail/Name/Port select statement.
#   - Active selection moves from C10 to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "            var result = _context.Services.AsQueryable().Where(t => t.Port > 10);`n            var res = result.ToList();"
$ws.Range("B6").Value = "SELECT s.Id, s.MaintainerEmail, s.Name, s.Port`n              FROM Service AS s`n              WHERE s.Port > 10"
$ws.Range("C6").Value = "            var result = _context.Services.AsEnumerable().Where(t => t.Port > 10);`n            var res = result.ToList();"
$ws.Range("D6").Value = "SELECT s.Id, s.MaintainerEmail, s.Name, s.Port`n              FROM Service AS s"

$ws.Range("C7").Value = "            var result = _context.Services.AsEnumerable().Select(t => t.Port);`n            var res = result.ToList();"
$ws.Range("D7").Value = "SELECT s.Id, s.MaintainerEmail, s.Name, s.Port`n              FROM Service AS s"

$ws.Range("A8").Value = "            var result = _context.Services.AsQueryable().Select(t => t.Port);`n            var res = result.Max();"
$ws.Range("C8").Value = "            var result = _context.Services.AsEnumerable().Select(t => t.Port);`n            var res = result.Max();"
$ws.Range("D8").Value = "SELECT s.Id, s.MaintainerEmail, s.Name, s.Port`n              FROM Service AS s"

$ws.Activate()
$ws.Range("C5").Select()
